$d = $word.ActiveDocument

$replacements = @(
    @("304÷4=", "314÷4="),
    @("844÷4=", "256÷6="),
    @("735÷7=", "452÷2="),
    @("569÷7=", "292÷5="),
    @("760÷4=", "258÷4="),
    @("499÷8=", "265÷3="),
    @("572÷5=", "155÷8="),
    @("749÷3=", "990÷5="),
    @("160÷3=", "134÷4="),
    @("599÷9=", "503÷4="),
    @("773÷9=", "458÷5="),
    @("628÷7=", "225÷4="),
    @("298÷9=", "289÷9="),
    @("708÷2=", "750÷5="),
    @("608÷7=", "714÷4="),
    @("519÷5=", "697÷3="),
    @("860÷5=", "103÷7="),
    @("534÷5=", "641÷4="),
    @("339÷2=", "591÷9="),
    @("480÷7=", "415÷4="),
    @("808÷9=", "670÷3="),
    @("130÷2=", "219÷2="),
    @("875÷6=", "417÷6="),
    @("666÷2=", "667÷3="),
    @("505÷7=", "595÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
